# Applies the "Updated symbol list" data refresh described by the diff:
# prices/volumes/hora bump for existing rows, plus a block of coin rows
# (9-20) that shift down by one slot with a new "One" entry inserted at
# the top and LiechtensteinCryptoassetsExchange rotating back in at row 20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''258.02'
$ws.Range("E2").Value = '''0.83%'
$ws.Range("G2").Value = '''4'
# Row 3
$ws.Range("D3").Value = '''27.12'
$ws.Range("E3").Value = '''-3.28%'
$ws.Range("G3").Value = '''4'
# Row 4
$ws.Range("D4").Value = '''4.916'
$ws.Range("E4").Value = '''-7.76%'
$ws.Range("G4").Value = '''4'
# Row 5
$ws.Range("D5").Value = '''0.05954'
$ws.Range("E5").Value = '''2.30%'
$ws.Range("G5").Value = '''4'
# Row 6
$ws.Range("D6").Value = '''6.689'
$ws.Range("E6").Value = '''-0.30%'
$ws.Range("G6").Value = '''4'
# Row 7
$ws.Range("D7").Value = '''0.8715'
$ws.Range("E7").Value = '''-0.19%'
$ws.Range("G7").Value = '''4'
# Row 8
$ws.Range("D8").Value = '''0.9586'
$ws.Range("E8").Value = '''6.99%'
$ws.Range("G8").Value = '''4'
# Row 9
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '''0.01062'
$ws.Range("E9").Value = '''1,659.34%'
$ws.Range("G9").Value = '''4'
# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1412'
$ws.Range("E10").Value = '''0.22%'
$ws.Range("G10").Value = '''4'
# Row 11
$ws.Range("D11").Value = '''0.07202'
$ws.Range("E11").Value = '''-0.23%'
$ws.Range("G11").Value = '''4'
# Row 12
$ws.Range("D12").Value = '''0.03132'
$ws.Range("E12").Value = '''-1.50%'
$ws.Range("G12").Value = '''4'
# Row 13
$ws.Range("D13").Value = '''0.09250'
$ws.Range("E13").Value = '''0.17%'
$ws.Range("G13").Value = '''4'
# Row 14
$ws.Range("D14").Value = '''0.001544'
$ws.Range("E14").Value = '''-1.70%'
$ws.Range("G14").Value = '''4'
# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005978'
$ws.Range("E15").Value = '''-1.12%'
$ws.Range("G15").Value = '''4'
# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.487'
$ws.Range("E16").Value = '''-0.30%'
$ws.Range("G16").Value = '''4'
# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.219'
$ws.Range("G17").Value = '''4'
# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.219'
$ws.Range("E18").Value = '''-2.32%'
$ws.Range("G18").Value = '''4'
# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3144'
$ws.Range("E19").Value = '''-0.71%'
$ws.Range("G19").Value = '''4'
# Row 20
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '''0.03564'
$ws.Range("E20").Value = '''4.16%'
$ws.Range("G20").Value = '''4'
# Row 21
$ws.Range("D21").Value = '''0.1306'
$ws.Range("E21").Value = '''-0.58%'
$ws.Range("G21").Value = '''4'
# Row 22
$ws.Range("E22").Value = '''0.20%'
$ws.Range("G22").Value = '''4'
# Row 23
$ws.Range("D23").Value = '''0.04257'
$ws.Range("E23").Value = '''2.35%'
$ws.Range("G23").Value = '''4'
# Row 24
$ws.Range("E24").Value = '''2.49%'
$ws.Range("G24").Value = '''4'
# Row 25
$ws.Range("D25").Value = '''0.001224'
$ws.Range("E25").Value = '''0.33%'
$ws.Range("G25").Value = '''4'
# Row 26
$ws.Range("D26").Value = '''0.004518'
$ws.Range("E26").Value = '''-7.22%'
$ws.Range("G26").Value = '''4'
# Row 27
$ws.Range("E27").Value = '''0.09%'
$ws.Range("G27").Value = '''4'
# Row 28
$ws.Range("D28").Value = '''0.0001493'
$ws.Range("E28").Value = '''-22.90%'
$ws.Range("G28").Value = '''4'
# Row 29
$ws.Range("G29").Value = '''4'
# Row 30
$ws.Range("G30").Value = '''4'
# Row 31
$ws.Range("G31").Value = '''4'
# Row 32
$ws.Range("G32").Value = '''4'
# Row 33
$ws.Range("G33").Value = '''4'
# Row 34
$ws.Range("G34").Value = '''4'
# Row 35
$ws.Range("G35").Value = '''4'
# Row 36
$ws.Range("G36").Value = '''4'
# Row 37
$ws.Range("G37").Value = '''4'
# Row 38
$ws.Range("G38").Value = '''4'
# Row 39
$ws.Range("G39").Value = '''4'
# Row 40
$ws.Range("D40").Value = '''0.03834'
$ws.Range("E40").Value = '''-0.75%'
$ws.Range("G40").Value = '''4'
# Row 41
$ws.Range("D41").Value = '''0.006579'
$ws.Range("E41").Value = '''14.32%'
$ws.Range("G41").Value = '''4'
# Row 42
$ws.Range("D42").Value = '''0.1101'
$ws.Range("E42").Value = '''0.09%'
$ws.Range("G42").Value = '''4'
# Row 43
$ws.Range("D43").Value = '''0.002200'
$ws.Range("E43").Value = '''-4.01%'
$ws.Range("G43").Value = '''4'
# Row 44
$ws.Range("D44").Value = '''0.01055'
$ws.Range("E44").Value = '''5.95%'
$ws.Range("G44").Value = '''4'
# Row 45
$ws.Range("D45").Value = '''0.00005488'
$ws.Range("E45").Value = '''4.12%'
$ws.Range("G45").Value = '''4'
# Row 46
$ws.Range("E46").Value = '''0.09%'
$ws.Range("G46").Value = '''4'
# Row 47
$ws.Range("E47").Value = '''28.45%'
$ws.Range("G47").Value = '''4'
# Row 48
$ws.Range("D48").Value = '''0.002126'
$ws.Range("E48").Value = '''-24.65%'
$ws.Range("G48").Value = '''4'
# Row 49
$ws.Range("E49").Value = '''0.09%'
$ws.Range("G49").Value = '''4'
# Row 50
$ws.Range("E50").Value = '''0.09%'
$ws.Range("G50").Value = '''4'
# Row 51
$ws.Range("G51").Value = '''4'
